$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.648.72"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "'3.444.79"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'592.25"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").Value = "'137.93"
$ws.Range("E6").Value = "  -5.27%  "
$ws.Range("D7").Value = "'3.444.21"
$ws.Range("E7").Value = "  -2.78%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.502"
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("D10").Value = "'7.34"
$ws.Range("E10").Value = "  -5.51%  "
$ws.Range("D11").Value = "'0.123"
$ws.Range("E11").Value = "  -7.68%  "
$ws.Range("E12").Value = "  -6.43%  "
$ws.Range("D13").Value = "'4.023.35"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("E14").Value = "  -8.97%  "
$ws.Range("D15").Value = "'26.55"
$ws.Range("E15").Value = "  -8.11%  "
$ws.Range("D16").Value = "'3.438.48"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("D17").Value = "'65.572.94"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "'9.86"
$ws.Range("E19").Value = "  -10.34%  "
$ws.Range("D20").Value = "'5.90"
$ws.Range("E20").Value = "  -4.68%  "
$ws.Range("D21").Value = "'13.77"
$ws.Range("E21").Value = "  -5.87%  "
$ws.Range("D22").Value = "'394.46"
$ws.Range("E22").Value = "  -5.26%  "
$ws.Range("D23").Value = "'0.555"
$ws.Range("E23").Value = "  -7.29%  "
$ws.Range("D24").Value = "'73.40"
$ws.Range("E24").Value = "  -5.60%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'3.585.08"
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("E27").Value = "  -7.17%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "'8.30"
$ws.Range("E29").Value = "  -8.17%  "
$ws.Range("D30").Value = "'7.20"
$ws.Range("E30").Value = "  -7.86%  "
$ws.Range("E31").Value = "  -8.65%  "
$ws.Range("D32").Value = "'3.449.25"
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -6.10%  "
$ws.Range("D35").Value = "'23.03"
$ws.Range("E35").Value = "  -5.62%  "
$ws.Range("D36").Value = "'172.68"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").Value = "'6.96"
$ws.Range("E37").Value = "  -7.73%  "
$ws.Range("D38").Value = "'1.19"
$ws.Range("E38").Value = "  -8.98%  "
$ws.Range("D39").Value = "'1.50"
$ws.Range("E39").Value = "  -6.47%  "
$ws.Range("E40").Value = "  -8.21%  "
$ws.Range("D41").Value = "'0.0770"
$ws.Range("E41").Value = "  -5.87%  "
$ws.Range("D42").Value = "'0.827"
$ws.Range("E42").Value = "  -3.71%  "
$ws.Range("D43").Value = "'43.79"
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'4.44"
$ws.Range("E45").Value = "  -12.33%  "
$ws.Range("D46").Value = "'1.63"
$ws.Range("E46").Value = "  -9.33%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.12"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'23.07"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("D49").Value = "'6.61"
$ws.Range("E49").Value = "  -6.49%  "
$ws.Range("D50").Value = "'2.11"
$ws.Range("E50").Value = "  -12.33%  "
$ws.Range("D51").Value = "'2.216.15"
$ws.Range("E51").Value = "  -6.55%  "
